$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Spring Questions")

# --- Sheet1: neutralize the "applyFont" style remnants on B64/B65/B66 ---
# (re-asserting the already-false Bold state nudges the style engine to
# drop the now-pointless font override, matching the source edit where
# these three cells lost their extra font-only style variant)
$ws1.Range("B64").Font.Bold = $false
$ws1.Range("B65").Font.Bold = $false
$ws1.Range("B66").Font.Bold = $false

# --- Spring Questions sheet: add a new answer row under "What is IOC" ---
$ws2.Rows.Item(6).EntireRow.Insert()
$ws2.Range("A6").Value = "IOC tells let the client class do not take object creation responsibility, let some external component create the objects required for the client and give it. This will give loose coupling."
$ws2.Range("A6").Font.Bold = $false
$ws2.Range("A6").WrapText = $true
$ws2.Rows.Item(6).RowHeight = 28.8

# --- Spring Questions sheet: restore the blank separator row before the
#     "What is AutoWiring" question (keeps the Q/A/blank rhythm) ---
$ws2.Rows.Item(12).EntireRow.Insert()

# --- Selections / active sheet, matching where the author ended up ---
[void]$ws1.Range("B3").Select()
[void]$ws2.Range("A17").Select()
$ws2.Activate()
